# "photo no pointer fo robo"
#
# Several backlog items on "Бэклог задач" are done (or no longer needed) and
# get moved down onto the "Бэклог Archive" sheet (marked O.k. / Not Required),
# one remaining backlog item gets flagged "!" (urgent), and a handful of new
# backlog ideas get appended to the end of the backlog list.

$wb = $excel.ActiveWorkbook
$backlog  = $wb.Worksheets.Item("Бэклог задач")
$archive  = $wb.Worksheets.Item("Бэклог Archive")

# Rows (in "Бэклог задач") that are finished and move to the archive, in the
# order they currently appear, together with the status they get archived
# under.
$rowsToArchive = @(17, 19, 20, 25, 28, 29, 31, 32)
$statusForRow = @{
    17 = "O.k.";
    19 = "O.k.";
    20 = "Not Required";
    25 = "O.k.";
    28 = "O.k.";
    29 = "O.k.";
    31 = "O.k.";
    32 = "O.k.";
}

# Find the first empty row at the bottom of the archive sheet.
$archiveRow = $archive.Cells.Item($archive.Rows.Count, 1).End(-4162).Row + 1

foreach ($r in $rowsToArchive) {
    $taskText = $backlog.Cells.Item($r, 2).Value2
    $archive.Cells.Item($archiveRow, 1).Value = $statusForRow[$r]
    $archive.Cells.Item($archiveRow, 2).Value = $taskText
    $archiveRow = $archiveRow + 1
}

# Now remove the archived rows from the backlog sheet, bottom-to-top so the
# row numbers of the rows still to be removed don't shift.
$sortedDesc = $rowsToArchive | Sort-Object -Descending
foreach ($r in $sortedDesc) {
    $backlog.Rows.Item($r).Delete()
}

# The remaining "Гостевой режим игры." row is flagged urgent.
$backlog.Cells.Item(19, 1).Value = "!"

# Append freshly captured backlog ideas to the end of the list.
$newIdeas = @(
    'Сделать очки  : "побед: 123" для каждой игры свой во время игры.',
    'сделать более яркую и понятную "пригласить"',
    "сделать надпись: свободные игроки",
    "сделать против рейтинга знак i в кружочке при наведении вывести кард инфо игрока.",
    "Сделать роботом еще более подстраиваемым",
    "сделать нотификацию"
)

$nextRow = $backlog.Cells.Item($backlog.Rows.Count, 2).End(-4162).Row + 1
foreach ($idea in $newIdeas) {
    $backlog.Cells.Item($nextRow, 2).Value = $idea
    $nextRow = $nextRow + 1
}
$lastIdeaRow = $nextRow - 1

# Reflect where each sheet was left selected: the archive view sits on the
# newly archived "photo no pointer fo robo" row, while the backlog view
# keeps focus on the freshly typed last idea - and stays the active tab.
$archive.Activate()
$archive.Cells.Item($archiveRow - 1, 1).Resize(1, 2).Select()

$backlog.Activate()
$backlog.Cells.Item($lastIdeaRow + 1, 2).Select()
